$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4019.5
$ws.Range("J51").Value = 4025
$ws.Range("L51").Value = 4025
$ws.Range("N51").Value = -4993

$ws.Range("H76").Value = 4102.3
$ws.Range("I76").Value = 3938.2144
$ws.Range("K76").Value = 3938.2144
$ws.Range("M76").Value = -3623.2144

$ws.Range("H79").Value = 4102.3
$ws.Range("I79").Value = 3938.2144
$ws.Range("K79").Value = 3938.2144
$ws.Range("M79").Value = -2846.2144

$ws.Range("H100").Value = 4064.5
$ws.Range("I100").Value = 1475.7142
$ws.Range("K100").Value = 1475.7142
$ws.Range("M100").Value = -934.7141999999999

$ws.Range("H138").Value = 3211.3726
$ws.Range("I138").Value = 2305.2104
$ws.Range("J138").Value = 3749.4062
$ws.Range("K138").Value = 6915.6312
$ws.Range("L138").Value = 11248.2186
$ws.Range("M138").Value = -1775.6312
$ws.Range("N138").Value = -21528.2186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 8999.5
$ws.Range("I29").Value = 13449
$ws.Range("J29").Value = 4550
$ws.Range("K29").Value = 13449
$ws.Range("L29").Value = 4550
$ws.Range("M29").Value = -13141
$ws.Range("N29").Value = -5166

$ws.Range("H54").Value = 24999.5
$ws.Range("J54").Value = 24999.5
$ws.Range("L54").Value = 24999.5
$ws.Range("N54").Value = -26537.5

$ws.Range("H122").Value = 2273.4
$ws.Range("I122").Value = 2273.4
$ws.Range("K122").Value = 6820.200000000001
$ws.Range("M122").Value = -4370.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3850.4443
$ws.Range("I20").Value = 3462.4614
$ws.Range("J20").Value = 4859.2
$ws.Range("K20").Value = 3462.4614
$ws.Range("L20").Value = 4859.2
$ws.Range("M20").Value = -3215.4614
$ws.Range("N20").Value = -5353.2

$ws.Range("H86").Value = 1437.5
$ws.Range("I86").Value = 1401.2
$ws.Range("J86").Value = 1498
$ws.Range("K86").Value = 1401.2
$ws.Range("L86").Value = 1498
$ws.Range("M86").Value = -278.2
$ws.Range("N86").Value = -3744

$ws.Range("H89").Value = 1437.5
$ws.Range("I89").Value = 1401.2
$ws.Range("J89").Value = 1498
$ws.Range("K89").Value = 7006
$ws.Range("L89").Value = 7490
$ws.Range("M89").Value = -1390
$ws.Range("N89").Value = -18722

$ws.Range("H105").Value = 3147.1482
$ws.Range("I105").Value = 2833.0557
$ws.Range("J105").Value = 3775.3333
$ws.Range("K105").Value = 2833.0557
$ws.Range("L105").Value = 3775.3333
$ws.Range("M105").Value = -1086.0557
$ws.Range("N105").Value = -7269.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7924.0425
$ws.Range("I31").Value = 2265.3572
$ws.Range("J31").Value = 10324.697
$ws.Range("K31").Value = 2265.3572
$ws.Range("L31").Value = 10324.697
$ws.Range("M31").Value = -1970.3572
$ws.Range("N31").Value = -10914.697

$ws.Range("H34").Value = 7924.0425
$ws.Range("I34").Value = 2265.3572
$ws.Range("J34").Value = 10324.697
$ws.Range("K34").Value = 2265.3572
$ws.Range("L34").Value = 10324.697
$ws.Range("M34").Value = -2063.3572
$ws.Range("N34").Value = -10728.697

$ws.Range("H111").Value = 72688.75
$ws.Range("J111").Value = 72688.75
$ws.Range("L111").Value = 72688.75
$ws.Range("N111").Value = -80868.75

$ws.Range("H122").Value = 1496.2307
$ws.Range("I122").Value = 1559.8572
$ws.Range("J122").Value = 1422
$ws.Range("K122").Value = 4679.571599999999
$ws.Range("L122").Value = 4266
$ws.Range("M122").Value = -2229.571599999999
$ws.Range("N122").Value = -9166

$ws.Range("H134").Value = 10524.789
$ws.Range("I134").Value = 10800.595
$ws.Range("J134").Value = 320
$ws.Range("K134").Value = 32401.785
$ws.Range("L134").Value = 960
$ws.Range("M134").Value = -29866.785
$ws.Range("N134").Value = -6030

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 659.625
$ws.Range("I31").Value = 659.625
$ws.Range("K31").Value = 1978.875
$ws.Range("M31").Value = -1690.875

$ws.Range("H34").Value = 1913.2858
$ws.Range("I34").Value = 848.5
$ws.Range("K34").Value = 2545.5
$ws.Range("M34").Value = -2461.5

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 3117
$ws.Range("I55").Value = 3117
$ws.Range("K55").Value = 9351
$ws.Range("M55").Value = -9174

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5955.091
$ws.Range("I70").Value = 6198.5
$ws.Range("J70").Value = 5901
$ws.Range("K70").Value = 6198.5
$ws.Range("L70").Value = 5901
$ws.Range("M70").Value = -5928.5
$ws.Range("N70").Value = -6441

$ws.Range("H73").Value = 5955.091
$ws.Range("I73").Value = 6198.5
$ws.Range("J73").Value = 5901
$ws.Range("K73").Value = 6198.5
$ws.Range("L73").Value = 5901
$ws.Range("M73").Value = -5262.5
$ws.Range("N73").Value = -7773

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 58433.61
$ws.Range("I22").Value = 144574
$ws.Range("K22").Value = 144574
$ws.Range("M22").Value = -144279

$ws.Range("H27").Value = 58433.61
$ws.Range("I27").Value = 144574
$ws.Range("K27").Value = 144574
$ws.Range("M27").Value = -144467

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()

$ws.Range("H61").Value = 14661.211
$ws.Range("I61").Value = 23572.908
$ws.Range("K61").Value = 23572.908
$ws.Range("M61").Value = -23370.908

$ws.Range("H68").Value = 5673.231
$ws.Range("I68").Value = 5392.2856
$ws.Range("J68").Value = 6001
$ws.Range("K68").Value = 5392.2856
$ws.Range("L68").Value = 6001
$ws.Range("M68").Value = -4643.2856
$ws.Range("N68").Value = -7499

$ws.Range("H71").Value = 5673.231
$ws.Range("I71").Value = 5392.2856
$ws.Range("J71").Value = 6001
$ws.Range("K71").Value = 26961.428
$ws.Range("L71").Value = 30005
$ws.Range("M71").Value = -23217.428
$ws.Range("N71").Value = -37493

$ws.Range("H100").Value = 9350.429
$ws.Range("I100").Value = 1598.75
$ws.Range("K100").Value = 1598.75
$ws.Range("M100").Value = -1057.75

$ws.Range("H113").Value = 14661.211
$ws.Range("I113").Value = 23572.908
$ws.Range("K113").Value = 23572.908
$ws.Range("M113").Value = -21402.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 250013250
$ws.Range("I2").Value = 17664.666
$ws.Range("K2").Value = 17664.666
$ws.Range("M2").Value = -17552.666

$ws.Range("H4").Value = 167599920
$ws.Range("I4").Value = 1850179.4
$ws.Range("J4").Value = 333349660
$ws.Range("K4").Value = 1850179.4
$ws.Range("L4").Value = 333349660
$ws.Range("M4").Value = -1850066.4
$ws.Range("N4").Value = -333349886
